# Updating unit mapping spreadsheet.
# - SW Units sheet: drop the "Too simple?" comment on the Ship/Player rows.
# - Functions sheet: record who (Melissa / Nick / ?) owns each function's status.

$wb = $excel.ActiveWorkbook

$wsUnits = $wb.Worksheets.Item("SW Units")
$wsFunctions = $wb.Worksheets.Item("Functions")

# --- SW Units sheet: remove the "Too simple?" comments from Ship & Player rows ---
$wsUnits.Range("C4").ClearContents()
$wsUnits.Range("C5").ClearContents()

# --- Functions sheet: add the new "Status" owner column (E) / completion column (F) ---

# Grid class block (rows 16-33) -> Melissa
$wsFunctions.Range("E16:E33").Value = "Melissa"

# Battleship class block (rows 34-44) -> ?
$wsFunctions.Range("E34:E44").Value = "?"

# Player class block (rows 45-48) -> Nick
$wsFunctions.Range("E45:E48").Value = "Nick"

# Randomizer class block (rows 49-58) -> ?
$wsFunctions.Range("E49:E58").Value = "?"

# Ship class block (rows 59-69) -> Nick
$wsFunctions.Range("E59:E69").Value = "Nick"

# Completed ("Done") items within the Ship block
$wsFunctions.Range("F59:F60").Value = "Done"
$wsFunctions.Range("F65:F68").Value = "Done"

# --- View state tweaks ---
$wsUnits.Activate()
$wsUnits.Range("A1").Select()

$wsFunctions.Activate()
$excel.ActiveWindow.Zoom = 85
